# Auto-generated edit script: updates market-price derived columns (H-N)
# across multiple sheets per the commit diff. Values are static data
# (no formulas in this workbook), so each changed cell is set directly.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 2498.6667
$ws.Cells.Item(18, 9).Value = 2163.3333
$ws.Cells.Item(18, 10).Value = 2834
$ws.Cells.Item(18, 11).Value = 2163.3333
$ws.Cells.Item(18, 12).Value = 2834
$ws.Cells.Item(18, 13).Value = -1879.3333
$ws.Cells.Item(18, 14).Value = -3402
$ws.Cells.Item(41, 8).Value = 1785.2307
$ws.Cells.Item(41, 9).Value = 2131.889
$ws.Cells.Item(41, 10).Value = 1005.25
$ws.Cells.Item(41, 11).Value = 2131.889
$ws.Cells.Item(41, 12).Value = 1005.25
$ws.Cells.Item(41, 13).Value = -1691.889
$ws.Cells.Item(41, 14).Value = -1885.25
$ws.Cells.Item(53, 8).Value = 1537.4
$ws.Cells.Item(53, 10).Value = 1668.5
$ws.Cells.Item(53, 12).Value = 1668.5
$ws.Cells.Item(53, 14).Value = -2942.5
$ws.Cells.Item(92, 8).Value = 5475.625
$ws.Cells.Item(92, 9).Value = 3054.3333
$ws.Cells.Item(92, 10).Value = 8588.714
$ws.Cells.Item(92, 11).Value = 3054.3333
$ws.Cells.Item(92, 12).Value = 8588.714
$ws.Cells.Item(92, 13).Value = -1806.3333
$ws.Cells.Item(92, 14).Value = -11084.714
$ws.Cells.Item(116, 8).Value = 20872.38
$ws.Cells.Item(116, 9).Value = 16132.9
$ws.Cells.Item(116, 11).Value = 16132.9
$ws.Cells.Item(116, 13).Value = -12690.9
$ws.Cells.Item(129, 8).Value = 6758.7646
$ws.Cells.Item(129, 9).Value = 866.6667
$ws.Cells.Item(129, 10).Value = 9972.637000000001
$ws.Cells.Item(129, 11).Value = 2600.0001
$ws.Cells.Item(129, 12).Value = 29917.911
$ws.Cells.Item(129, 13).Value = 2399.9999
$ws.Cells.Item(129, 14).Value = -39917.911
$ws.Cells.Item(131, 8).Value = 1895724.2
$ws.Cells.Item(131, 9).Value = 1942.9
$ws.Cells.Item(131, 11).Value = 5828.700000000001
$ws.Cells.Item(131, 13).Value = -788.7000000000007
$ws.Cells.Item(132, 8).Value = 2369.087
$ws.Cells.Item(132, 9).Value = 2333.0476
$ws.Cells.Item(132, 10).Value = 2747.5
$ws.Cells.Item(132, 11).Value = 6999.1428
$ws.Cells.Item(132, 12).Value = 8242.5
$ws.Cells.Item(132, 13).Value = -4469.1428
$ws.Cells.Item(132, 14).Value = -13302.5
$ws.Cells.Item(138, 8).Value = 5145.0586
$ws.Cells.Item(138, 10).Value = 6438.48
$ws.Cells.Item(138, 12).Value = 19315.44
$ws.Cells.Item(138, 14).Value = -29595.44

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6796.393
$ws.Cells.Item(32, 9).Value = 5336.519
$ws.Cells.Item(32, 11).Value = 5336.519
$ws.Cells.Item(32, 13).Value = -5049.519
$ws.Cells.Item(45, 8).Value = 2804456.8
$ws.Cells.Item(45, 9).Value = 7938180
$ws.Cells.Item(45, 11).Value = 7938180
$ws.Cells.Item(45, 13).Value = -7937803
$ws.Cells.Item(60, 8).Value = 90276.55499999999
$ws.Cells.Item(60, 9).Value = 90276.55499999999
$ws.Cells.Item(60, 11).Value = 90276.55499999999
$ws.Cells.Item(60, 13).Value = -89543.55499999999
$ws.Cells.Item(74, 8).Value = 2734.6667
$ws.Cells.Item(74, 10).Value = 2857
$ws.Cells.Item(74, 12).Value = 2857
$ws.Cells.Item(74, 14).Value = -4605
$ws.Cells.Item(77, 8).Value = 2734.6667
$ws.Cells.Item(77, 10).Value = 2857
$ws.Cells.Item(77, 12).Value = 14285
$ws.Cells.Item(77, 14).Value = -23021
$ws.Cells.Item(102, 8).Value = 5347.647
$ws.Cells.Item(102, 9).Value = 3856.6667
$ws.Cells.Item(102, 11).Value = 3856.6667
$ws.Cells.Item(102, 13).Value = -2234.6667
$ws.Cells.Item(122, 8).Value = 10833
$ws.Cells.Item(122, 9).Value = 15999
$ws.Cells.Item(122, 11).Value = 47997
$ws.Cells.Item(122, 13).Value = -45547
$ws.Cells.Item(132, 8).Value = 3126819.5
$ws.Cells.Item(132, 9).Value = 1877.4828
$ws.Cells.Item(132, 11).Value = 5632.4484
$ws.Cells.Item(132, 13).Value = -3102.4484

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 13957.6
$ws.Cells.Item(22, 9).Value = 930
$ws.Cells.Item(22, 10).Value = 33499
$ws.Cells.Item(22, 11).Value = 930
$ws.Cells.Item(22, 12).Value = 33499
$ws.Cells.Item(22, 13).Value = -757
$ws.Cells.Item(22, 14).Value = -33845
$ws.Cells.Item(94, 8).Value = 1927.7142
$ws.Cells.Item(94, 9).Value = 2683.3333
$ws.Cells.Item(94, 10).Value = 920.2222
$ws.Cells.Item(94, 11).Value = 2683.3333
$ws.Cells.Item(94, 12).Value = 920.2222
$ws.Cells.Item(94, 13).Value = -2232.3333
$ws.Cells.Item(94, 14).Value = -1822.2222
$ws.Cells.Item(134, 8).Value = 4001952
$ws.Cells.Item(134, 9).Value = 2060.7827
$ws.Cells.Item(134, 10).Value = 50000700
$ws.Cells.Item(134, 11).Value = 6182.348100000001
$ws.Cells.Item(134, 12).Value = 150002100
$ws.Cells.Item(134, 13).Value = -3647.348100000001
$ws.Cells.Item(134, 14).Value = -150007170

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 7693351
$ws.Cells.Item(16, 9).Value = 11111709
$ws.Cells.Item(16, 11).Value = 11111709
$ws.Cells.Item(16, 13).Value = -11111422
$ws.Cells.Item(31, 8).Value = 38464576
$ws.Cells.Item(31, 9).Value = 47621450
$ws.Cells.Item(31, 10).Value = 5702.8
$ws.Cells.Item(31, 11).Value = 47621450
$ws.Cells.Item(31, 12).Value = 5702.8
$ws.Cells.Item(31, 13).Value = -47621155
$ws.Cells.Item(31, 14).Value = -6292.8
$ws.Cells.Item(34, 8).Value = 38464576
$ws.Cells.Item(34, 9).Value = 47621450
$ws.Cells.Item(34, 10).Value = 5702.8
$ws.Cells.Item(34, 11).Value = 47621450
$ws.Cells.Item(34, 12).Value = 5702.8
$ws.Cells.Item(34, 13).Value = -47621248
$ws.Cells.Item(34, 14).Value = -6106.8
$ws.Cells.Item(86, 8).Value = 6212.923
$ws.Cells.Item(86, 9).Value = 5730.3335
$ws.Cells.Item(86, 11).Value = 5730.3335
$ws.Cells.Item(86, 13).Value = -4607.3335
$ws.Cells.Item(89, 8).Value = 6212.923
$ws.Cells.Item(89, 9).Value = 5730.3335
$ws.Cells.Item(89, 11).Value = 28651.6675
$ws.Cells.Item(89, 13).Value = -23035.6675
$ws.Cells.Item(113, 8).Value = 7693351
$ws.Cells.Item(113, 9).Value = 11111709
$ws.Cells.Item(113, 11).Value = 11111709
$ws.Cells.Item(113, 13).Value = -11109539
$ws.Cells.Item(122, 8).Value = 2790.3635
$ws.Cells.Item(122, 9).Value = 2991.875
$ws.Cells.Item(122, 11).Value = 8975.625
$ws.Cells.Item(122, 13).Value = -6525.625
$ws.Cells.Item(132, 8).Value = 1578.3871
$ws.Cells.Item(132, 9).Value = 1607.75
$ws.Cells.Item(132, 11).Value = 4823.25
$ws.Cells.Item(132, 13).Value = -2293.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 952.625
$ws.Cells.Item(5, 9).Value = 972.5
$ws.Cells.Item(5, 10).Value = 932.75
$ws.Cells.Item(5, 11).Value = 2917.5
$ws.Cells.Item(5, 12).Value = 2798.25
$ws.Cells.Item(5, 13).Value = -2805.5
$ws.Cells.Item(5, 14).Value = -3022.25
$ws.Cells.Item(46, 8).Value = 17066
$ws.Cells.Item(46, 10).Value = 17066
$ws.Cells.Item(46, 12).Value = 51198
$ws.Cells.Item(46, 14).Value = -51380
$ws.Cells.Item(56, 8).Value = 8607.182000000001
$ws.Cells.Item(56, 9).Value = 8607.182000000001
$ws.Cells.Item(56, 11).Value = 8607.182000000001
$ws.Cells.Item(56, 13).Value = -8077.182000000001
$ws.Cells.Item(121, 8).Value = 4274.2
$ws.Cells.Item(121, 9).Value = 497.5
$ws.Cells.Item(121, 11).Value = 1492.5
$ws.Cells.Item(121, 13).Value = -182.5
$ws.Cells.Item(122, 8).Value = 110733.336
$ws.Cells.Item(122, 9).Value = 165600
$ws.Cells.Item(122, 10).Value = 1000
$ws.Cells.Item(122, 11).Value = 1490400
$ws.Cells.Item(122, 12).Value = 9000
$ws.Cells.Item(122, 13).Value = -1487950
$ws.Cells.Item(122, 14).Value = -13900
$ws.Cells.Item(135, 8).Value = 952.625
$ws.Cells.Item(135, 9).Value = 972.5
$ws.Cells.Item(135, 10).Value = 932.75
$ws.Cells.Item(135, 11).Value = 8752.5
$ws.Cells.Item(135, 12).Value = 8394.75
$ws.Cells.Item(135, 13).Value = -6217.5
$ws.Cells.Item(135, 14).Value = -13464.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 1857024
$ws.Cells.Item(113, 9).Value = 4761.6
$ws.Cells.Item(113, 10).Value = 3709286.5
$ws.Cells.Item(113, 11).Value = 4761.6
$ws.Cells.Item(113, 12).Value = 3709286.5
$ws.Cells.Item(113, 13).Value = -2591.6
$ws.Cells.Item(113, 14).Value = -3713626.5
$ws.Cells.Item(132, 8).Value = 2568255.8
$ws.Cells.Item(132, 9).Value = 4421.0625
$ws.Cells.Item(132, 11).Value = 13263.1875
$ws.Cells.Item(132, 13).Value = -10733.1875

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 4874
$ws.Cells.Item(22, 9).Value = 4832.3335
$ws.Cells.Item(22, 11).Value = 4832.3335
$ws.Cells.Item(22, 13).Value = -4537.3335
$ws.Cells.Item(27, 8).Value = 4874
$ws.Cells.Item(27, 9).Value = 4832.3335
$ws.Cells.Item(27, 11).Value = 4832.3335
$ws.Cells.Item(27, 13).Value = -4725.3335
$ws.Cells.Item(42, 8).Value = 28805
$ws.Cells.Item(42, 9).Value = 28805
$ws.Cells.Item(42, 11).Value = 28805
$ws.Cells.Item(42, 13).Value = -28242
$ws.Cells.Item(49, 8).Value = 28805
$ws.Cells.Item(49, 9).Value = 28805
$ws.Cells.Item(49, 11).Value = 28805
$ws.Cells.Item(49, 13).Value = -28658
$ws.Cells.Item(61, 8).Value = 55559916
$ws.Cells.Item(61, 9).Value = 71432540
$ws.Cells.Item(61, 11).Value = 71432540
$ws.Cells.Item(61, 13).Value = -71432338
$ws.Cells.Item(113, 8).Value = 55559916
$ws.Cells.Item(113, 9).Value = 71432540
$ws.Cells.Item(113, 11).Value = 71432540
$ws.Cells.Item(113, 13).Value = -71430370
$ws.Cells.Item(122, 8).Value = 3982.36
$ws.Cells.Item(122, 9).Value = 3287.2683
$ws.Cells.Item(122, 11).Value = 9861.804900000001
$ws.Cells.Item(122, 13).Value = -7411.804900000001
$ws.Cells.Item(132, 8).Value = 4532.1816
$ws.Cells.Item(132, 9).Value = 2790
$ws.Cells.Item(132, 10).Value = 6622.8
$ws.Cells.Item(132, 11).Value = 8370
$ws.Cells.Item(132, 12).Value = 19868.4
$ws.Cells.Item(132, 13).Value = -5840
$ws.Cells.Item(132, 14).Value = -24928.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(87, 8).Value = 99999
$ws.Cells.Item(87, 10).Value = 99999
$ws.Cells.Item(87, 12).Value = 99999
$ws.Cells.Item(87, 14).Value = -102495
$ws.Cells.Item(90, 8).Value = 99999
$ws.Cells.Item(90, 10).Value = 99999
$ws.Cells.Item(90, 12).Value = 299997
$ws.Cells.Item(90, 14).Value = -312477
$ws.Cells.Item(132, 8).Value = 279493.9
$ws.Cells.Item(132, 9).Value = 1693.3846
$ws.Cells.Item(132, 11).Value = 5080.1538
$ws.Cells.Item(132, 13).Value = -2550.1538

